$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Repayment schedule")

# Insert a new blank column before column N (14th column), shifting the
# existing "Late", "heading" and "Outstanding" columns one place to the right.
$ws.Columns.Item(14).Insert()

# The freshly inserted column inherits Excel's default "insert column"
# width behaviour; set it explicitly to match (stored width 11).
$ws.Columns.Item(14).ColumnWidth = 10.166666666666666

# Make "Repayment schedule" the active sheet (previously "Transactions" was active)
$ws.Activate()

# Update the selected cell on the Repayment schedule sheet
$ws.Range("R9").Select()
